$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.949.43"
$ws.Range("E2").Value = "  +4.61%  "
$ws.Range("D3").Value = "2.233.79"
$ws.Range("E3").Value = "  +3.59%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "259.89"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "81.83"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +12.14%  "
$ws.Range("E7").Value = "  +3.35%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("E9").Value = "  +3.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.66"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +10.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0929"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.06"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +5.19%  "
$ws.Range("E13").Value = "  +2.87%  "
$ws.Range("D14").Value = "2.565.03"
$ws.Range("E14").Value = "  +3.13%  "
$ws.Range("E15").Value = "  +3.72%  "
$ws.Range("D16").Value = "2.239.64"
$ws.Range("E17").Value = "  +2.97%  "
$ws.Range("D18").Value = "43.849.45"
$ws.Range("E18").Value = "  +4.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.20"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("E21").Value = "  +4.25%  "
$ws.Range("E22").Value = "  +10.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.88"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.32"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.93%  "
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("E26").Value = "  +3.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "40.88"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +11.49%  "
$ws.Range("E28").Value = "  +1.65%  "
$ws.Range("E29").Value = "  +2.97%  "
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.80"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0904"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +14.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.63"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.35"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +5.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.116"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +9.19%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0372"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +13.92%  "
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.123"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.57"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +7.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.06"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +9.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.02"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +28.18%  "
$ws.Range("E41").Value = "  +4.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "63.44"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +8.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.55"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +8.41%  "
$ws.Range("E44").Value = "  +4.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.40"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.46"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.87%  "
$ws.Range("E47").Value = "  +2.53%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.57"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +30.36%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.13"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.445"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.10%  "
$ws.Range("E51").Value = "  +3.54%  "
